$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each "round" block spans 16 rows (one header/contribution row per player +
# the totals row), starting at row 2 and stepping by 16, for rounds 1..10.
# The A column was merged across each block showing the round number only
# once; unmerge it, fill the round number down every row of the block, and
# re-apply the "top" cell's format (style index 9) to the rest of the block
# so every cell in the column carries that same style instead of the
# merged-cell border variants.
for ($round = 1; $round -le 10; $round++) {
    $startRow = 2 + ($round - 1) * 16
    $endRow = $startRow + 15

    $blockRange = $ws.Range("A$startRow`:A$endRow")
    $blockRange.UnMerge()
    $blockRange.Value = $round

    $topCell = $ws.Range("A$startRow")
    $restRange = $ws.Range("A$($startRow + 1)`:A$endRow")
    $topCell.Copy()
    $restRange.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Restore the view state captured at save time: no pinned top-left cell and
# the active selection sitting on E156.
$ws.Range("E156").Select()
